$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# The "Periodo Mora" (E16:E27) and "Valor Mora" (F16:F27) blocks are
# reversed: the oldest-period row (16) and newest-period row (27) swap
# places, cascading down the whole 12-row block.

$periods = @("2104","2103","2101","2012","2011","2010","2009","2008","2007","2006","2005","2004")
$valorMora = @(41600,48000,48000,48000,48000,48000,48000,48000,48000,48000,48000,48000)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valorMora[$i]
}
